$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.23%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.39%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.122"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.61%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08108"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.92%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.139"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.01%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.12%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.152"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.07%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9264"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.83%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.84%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1870"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.64%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09189"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.91%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03604"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09903"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001436"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.13%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005729"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.03%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.73%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'3.74%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3367"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.89%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1332"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.02%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.120"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.51%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2221"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.73%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04566"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.17%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004704"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.95%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001252"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-21.92%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004509"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-5.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01954"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.27%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1391"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.56%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007827"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.50%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.18%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01165"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.01%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006512"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.87%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.34%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'39.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-17.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001703"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-14.79%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.34%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.34%"
$ws.Range("E51").Style = "Normal"
